# Applies the per-cell numeric updates captured in the commit diff for
# Sheets/Ixion_Profits.xlsx, scoped to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# worksheets of this workbook.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 73188.92999999999
$ws.Range("I111").Value = 1806.8572
$ws.Range("J111").Value = 144571
$ws.Range("K111").Value = 5420.571599999999
$ws.Range("L111").Value = 433713
$ws.Range("M111").Value = -2353.571599999999
$ws.Range("N111").Value = -439847

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 7000
$ws.Range("H32").Value = 4588.3564
$ws.Range("I32").Value = 2849.6758
$ws.Range("J32").Value = 14485.462
$ws.Range("K32").Value = 2849.6758
$ws.Range("L32").Value = 14485.462
$ws.Range("M32").Value = -2562.6758
$ws.Range("N32").Value = -15059.462
$ws.Range("H63").Value = 90911640
$ws.Range("I63").Value = 100002500
$ws.Range("K63").Value = 100002500
$ws.Range("M63").Value = -100001814
$ws.Range("H66").Value = 90911640
$ws.Range("I66").Value = 100002500
$ws.Range("K66").Value = 500012500
$ws.Range("M66").Value = -500009068
$ws.Range("H74").Value = 1602.1052
$ws.Range("I74").Value = 1549.3334
$ws.Range("K74").Value = 1549.3334
$ws.Range("M74").Value = -675.3334
$ws.Range("H77").Value = 1602.1052
$ws.Range("I77").Value = 1549.3334
$ws.Range("K77").Value = 7746.666999999999
$ws.Range("M77").Value = -3378.666999999999
$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5808
$ws.Range("H109").Value = 50188.5
$ws.Range("J109").Value = 50188.5
$ws.Range("L109").Value = 50188.5
$ws.Range("N109").Value = -52962.5
$ws.Range("H122").Value = 1352355.6
$ws.Range("I122").Value = 5135285
$ws.Range("J122").Value = 1309.5
$ws.Range("K122").Value = 15405855
$ws.Range("L122").Value = 3928.5
$ws.Range("M122").Value = -15403405
$ws.Range("N122").Value = -8828.5
$ws.Range("H132").Value = 2865.4634
$ws.Range("I132").Value = 1661.7059
$ws.Range("J132").Value = 8712.286
$ws.Range("K132").Value = 4985.1177
$ws.Range("L132").Value = 26136.858
$ws.Range("M132").Value = -2455.1177
$ws.Range("N132").Value = -31196.858

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20886.182
$ws.Range("I82").Value = 2402.1667
$ws.Range("J82").Value = 27817.688
$ws.Range("K82").Value = 2402.1667
$ws.Range("L82").Value = 27817.688
$ws.Range("M82").Value = -2019.1667
$ws.Range("N82").Value = -28583.688
$ws.Range("H85").Value = 20886.182
$ws.Range("I85").Value = 2402.1667
$ws.Range("J85").Value = 27817.688
$ws.Range("K85").Value = 2402.1667
$ws.Range("L85").Value = 27817.688
$ws.Range("M85").Value = -1076.1667
$ws.Range("N85").Value = -30469.688
$ws.Range("H86").Value = 9261245
$ws.Range("I86").Value = 12822353
$ws.Range("J86").Value = 2362.1
$ws.Range("K86").Value = 12822353
$ws.Range("L86").Value = 2362.1
$ws.Range("M86").Value = -12821230
$ws.Range("N86").Value = -4608.1
$ws.Range("H89").Value = 9261245
$ws.Range("I89").Value = 12822353
$ws.Range("J89").Value = 2362.1
$ws.Range("K89").Value = 64111765
$ws.Range("L89").Value = 11810.5
$ws.Range("M89").Value = -64106149
$ws.Range("N89").Value = -23042.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2404935.2
$ws.Range("I16").Value = 5129093
$ws.Range("J16").Value = 1266.6471
$ws.Range("K16").Value = 5129093
$ws.Range("L16").Value = 1266.6471
$ws.Range("M16").Value = -5128806
$ws.Range("N16").Value = -1840.6471
$ws.Range("H31").Value = 2703.1667
$ws.Range("I31").Value = 1356.814
$ws.Range("K31").Value = 1356.814
$ws.Range("M31").Value = -1061.814
$ws.Range("H34").Value = 2703.1667
$ws.Range("I34").Value = 1356.814
$ws.Range("K34").Value = 1356.814
$ws.Range("M34").Value = -1154.814
$ws.Range("H41").Value = 15000
$ws.Range("J41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15856
$ws.Range("H45").Value = 4750
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1407
$ws.Range("H50").Value = 21493.334
$ws.Range("J50").Value = 21493.334
$ws.Range("L50").Value = 21493.334
$ws.Range("N50").Value = -22743.334
$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26472
$ws.Range("H59").Value = 23398.715
$ws.Range("J59").Value = 23398.715
$ws.Range("L59").Value = 23398.715
$ws.Range("N59").Value = -25688.715
$ws.Range("H60").Value = 25000
$ws.Range("J60").Value = 25000
$ws.Range("L60").Value = 25000
$ws.Range("N60").Value = -26022
$ws.Range("H61").Value = 25000
$ws.Range("J61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25696
$ws.Range("H68").Value = 29533
$ws.Range("J68").Value = 29533
$ws.Range("L68").Value = 29533
$ws.Range("N68").Value = -31031
$ws.Range("H71").Value = 29533
$ws.Range("J71").Value = 29533
$ws.Range("L71").Value = 88599
$ws.Range("N71").Value = -96087
$ws.Range("H74").Value = 22919.8
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 22919.8
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 22919.8
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -24667.8
$ws.Range("H77").Value = 22919.8
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 22919.8
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 68759.39999999999
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -77495.39999999999
$ws.Range("H113").Value = 2404935.2
$ws.Range("I113").Value = 5129093
$ws.Range("J113").Value = 1266.6471
$ws.Range("K113").Value = 5129093
$ws.Range("L113").Value = 1266.6471
$ws.Range("M113").Value = -5126923
$ws.Range("N113").Value = -5606.6471

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 250676.38
$ws.Range("I5").Value = 746.9167
$ws.Range("J5").Value = 500605.84
$ws.Range("K5").Value = 2240.7501
$ws.Range("L5").Value = 1501817.52
$ws.Range("M5").Value = -2128.7501
$ws.Range("N5").Value = -1502041.52
$ws.Range("H86").Value = 650.375
$ws.Range("I86").Value = 853.3333
$ws.Range("K86").Value = 2559.9999
$ws.Range("M86").Value = -1373.9999
$ws.Range("H89").Value = 650.375
$ws.Range("I89").Value = 853.3333
$ws.Range("K89").Value = 7679.9997
$ws.Range("M89").Value = -1751.9997
$ws.Range("H131").Value = 3334346
$ws.Range("I131").Value = 16667105
$ws.Range("J131").Value = 1156.25
$ws.Range("K131").Value = 50001315
$ws.Range("L131").Value = 3468.75
$ws.Range("M131").Value = -49996275
$ws.Range("N131").Value = -13548.75
$ws.Range("H135").Value = 250676.38
$ws.Range("I135").Value = 746.9167
$ws.Range("J135").Value = 500605.84
$ws.Range("K135").Value = 6722.2503
$ws.Range("L135").Value = 4505452.560000001
$ws.Range("M135").Value = -4187.2503
$ws.Range("N135").Value = -4510522.560000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6322.4346
$ws.Range("I70").Value = 6244.3335
$ws.Range("J70").Value = 6603.6
$ws.Range("K70").Value = 6244.3335
$ws.Range("L70").Value = 6603.6
$ws.Range("M70").Value = -5974.3335
$ws.Range("N70").Value = -7143.6
$ws.Range("H73").Value = 6322.4346
$ws.Range("I73").Value = 6244.3335
$ws.Range("J73").Value = 6603.6
$ws.Range("K73").Value = 6244.3335
$ws.Range("L73").Value = 6603.6
$ws.Range("M73").Value = -5308.3335
$ws.Range("N73").Value = -8475.6
$ws.Range("H113").Value = 43479452
$ws.Range("I113").Value = 83334264
$ws.Range("K113").Value = 83334264
$ws.Range("M113").Value = -83332094
$ws.Range("H126").Value = 7860.7646
$ws.Range("I126").Value = 10177.667
$ws.Range("K126").Value = 30533.001
$ws.Range("M126").Value = -28063.001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4446227.5
$ws.Range("I22").Value = 18519182
$ws.Range("J22").Value = 2136.2104
$ws.Range("K22").Value = 18519182
$ws.Range("L22").Value = 2136.2104
$ws.Range("M22").Value = -18518887
$ws.Range("N22").Value = -2726.2104
$ws.Range("H27").Value = 4446227.5
$ws.Range("I27").Value = 18519182
$ws.Range("J27").Value = 2136.2104
$ws.Range("K27").Value = 18519182
$ws.Range("L27").Value = 2136.2104
$ws.Range("M27").Value = -18519075
$ws.Range("N27").Value = -2350.2104
$ws.Range("H46").Value = 27778992
$ws.Range("I46").Value = 47619972
$ws.Range("J46").Value = 1620
$ws.Range("K46").Value = 47619972
$ws.Range("L46").Value = 1620
$ws.Range("M46").Value = -47619784
$ws.Range("N46").Value = -1996
$ws.Range("H61").Value = 3796.1538
$ws.Range("I61").Value = 2924.1667
$ws.Range("J61").Value = 4543.5713
$ws.Range("K61").Value = 2924.1667
$ws.Range("L61").Value = 4543.5713
$ws.Range("M61").Value = -2722.1667
$ws.Range("N61").Value = -4947.5713
$ws.Range("H113").Value = 3796.1538
$ws.Range("I113").Value = 2924.1667
$ws.Range("J113").Value = 4543.5713
$ws.Range("K113").Value = 2924.1667
$ws.Range("L113").Value = 4543.5713
$ws.Range("M113").Value = -754.1667000000002
$ws.Range("N113").Value = -8883.5713
$ws.Range("H132").Value = 28649546
$ws.Range("I132").Value = 45836676
$ws.Range("J132").Value = 4331
$ws.Range("K132").Value = 137510028
$ws.Range("L132").Value = 12993
$ws.Range("M132").Value = -137507498
$ws.Range("N132").Value = -18053

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 698.2105
$ws.Range("I126").Value = 544.4706
$ws.Range("K126").Value = 1633.4118
$ws.Range("M126").Value = 836.5882000000001
$ws.Range("H132").Value = 1264.683
$ws.Range("I132").Value = 973.3
$ws.Range("J132").Value = 2059.3635
$ws.Range("K132").Value = 2919.9
$ws.Range("L132").Value = 6178.0905
$ws.Range("M132").Value = -389.8999999999996
$ws.Range("N132").Value = -11238.0905
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
